$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "missing ] symbol" bug: several of the roi (bounding-box) strings in
# column A had picked up a stray extra closing bracket at the very end
# (e.g. "...]]]]" instead of the correct "...]]]"). Strip that one extra
# character wherever it shows up.
for ($r = 1; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()
    if ($v -ne $null -and $v.EndsWith("]]]]")) {
        $cell.Value = $v.Substring(0, $v.Length - 1)
    }
}

# Update the on-screen selection / scroll position that Excel persisted in
# the sheet view (previously frozen at top-left cell A8 with B18 selected).
$ws.Range("E4").Select()
